$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(4)
$chart = $sh.Chart
Write-Host "Chart Name: $($chart.Name)"
Write-Host "HasTitle: $($chart.HasTitle)"
